$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U2").Formula = "1.92"
$ws.Range("V2").Formula = "1.77"
$ws.Range("M3").Formula = "1.1"
$ws.Range("O3").Formula = "1.44"
$ws.Range("P3").Formula = "2.75"
$ws.Range("U3").Formula = "1.92"
$ws.Range("V3").Formula = "1.77"
$ws.Range("M4").Formula = "1.08"
$ws.Range("O4").Formula = "1.44"
$ws.Range("P4").Formula = "2.75"
$ws.Range("V4").Formula = "1.63"
$ws.Range("M5").Formula = "1.11"
$ws.Range("O5").Formula = "1.5"
$ws.Range("Q5").Formula = "2.5"
$ws.Range("R5").Formula = "1.5"
$ws.Range("M6").Formula = "1.08"
$ws.Range("O6").Formula = "1.4"
$ws.Range("R6").Formula = "1.6"
$ws.Range("M7").Formula = "1.05"
$ws.Range("O7").Formula = "1.29"
$ws.Range("M8").Formula = "1.05"
$ws.Range("O8").Formula = "1.29"
$ws.Range("I9").Formula = "1.83"
$ws.Range("G11").Formula = "1.7"
$ws.Range("M11").Formula = "1.03"
$ws.Range("O11").Formula = "1.18"
$ws.Range("J12").Formula = "2.88"
$ws.Range("AA13").Formula = "29"
$ws.Range("AB13").Formula = "29"
$ws.Range("AE13").Formula = "13"
$ws.Range("AG13").Formula = "101"
$ws.Range("AH13").Formula = "12"
$ws.Range("AI13").Formula = "11"
$ws.Range("AN13").Formula = "7"
$ws.Range("AQ13").Formula = "67"
$ws.Range("AR13").Formula = "67"
$ws.Range("AX13").Formula = "8"
$ws.Range("AY13").Formula = "13"
$ws.Range("AZ13").Formula = "21"
$ws.Range("G13").Formula = "4.33"
$ws.Range("H13").Formula = "4.5"
$ws.Range("I13").Formula = "1.65"
$ws.Range("J13").Formula = "4.33"
$ws.Range("K13").Formula = "2.63"
$ws.Range("L13").Formula = "2.1"
$ws.Range("M13").Formula = "1.02"
$ws.Range("N13").Formula = "21"
$ws.Range("W13").Formula = "21"
$ws.Range("X13").Formula = "29"
$ws.Range("Z13").Formula = "51"
$ws.Range("Q14").Formula = "1.37"
$ws.Range("R14").Formula = "2.87"
$ws.Range("K15").Formula = "2.38"
$ws.Range("Q15").Formula = "1.63"
$ws.Range("AJ16").Formula = "17"
$ws.Range("AX16").Formula = "29"
$ws.Range("AZ16").Formula = "101"
$ws.Range("G16").Formula = "1.53"
$ws.Range("H16").Formula = "4.2"
$ws.Range("I16").Formula = "5.5"
$ws.Range("J16").Formula = "2.1"
$ws.Range("Q16").Formula = "1.85"
$ws.Range("R16").Formula = "2"
$ws.Range("U16").Formula = "1.83"
$ws.Range("V16").Formula = "1.83"
$ws.Range("X16").Formula = "7.5"
$ws.Range("M17").Formula = "1.02"
$ws.Range("O17").Formula = "1.13"
$ws.Range("M18").Formula = "1.05"
$ws.Range("O18").Formula = "1.29"
$ws.Range("M19").Formula = "1.03"
$ws.Range("O19").Formula = "1.22"
$ws.Range("AB20").Formula = "41"
$ws.Range("AE20").Formula = "17"
$ws.Range("AI20").Formula = "8.5"
$ws.Range("AK20").Formula = "11"
$ws.Range("AN20").Formula = "8"
$ws.Range("AO20").Formula = "29"
$ws.Range("AQ20").Formula = "101"
$ws.Range("AR20").Formula = "101"
$ws.Range("AU20").Formula = "8"
$ws.Range("AX20").Formula = "7"
$ws.Range("AZ20").Formula = "19"
$ws.Range("G20").Formula = "5.75"
$ws.Range("H20").Formula = "4.75"
$ws.Range("I20").Formula = "1.44"
$ws.Range("L20").Formula = "1.91"
$ws.Range("M20").Formula = "1.02"
$ws.Range("O20").Formula = "1.14"
$ws.Range("Q20").Formula = "1.48"
$ws.Range("R20").Formula = "2.6"
$ws.Range("U20").Formula = "1.63"
$ws.Range("V20").Formula = "2.1"
$ws.Range("W20").Formula = "21"
$ws.Range("Z20").Formula = "67"
$ws.Range("O21").Formula = "1.1"
$ws.Range("Q21").Formula = "1.33"
$ws.Range("U21").Formula = "1.5"
$ws.Range("V21").Formula = "2.37"
$ws.Range("M22").Formula = "1.04"
$ws.Range("O22").Formula = "1.22"
$ws.Range("Q22").Formula = "1.7"
$ws.Range("U22").Formula = "1.54"
$ws.Range("V23").Formula = "1.69"
$ws.Range("V24").Formula = "1.69"
$ws.Range("Q25").Formula = "1.88"
$ws.Range("R25").Formula = "1.98"
$ws.Range("U25").Formula = "1.77"
$ws.Range("V25").Formula = "1.92"
$ws.Range("U26").Formula = "1.58"
$ws.Range("G27").Formula = "2.1"
$ws.Range("I27").Formula = "3.75"
$ws.Range("J27").Formula = "2.88"
$ws.Range("L27").Formula = "4.33"
$ws.Range("Q28").Formula = "1.98"
$ws.Range("R28").Formula = "1.88"
$ws.Range("U30").Formula = "1.77"
$ws.Range("V30").Formula = "1.87"
$ws.Range("G32").Formula = "2.25"
$ws.Range("U32").Formula = "1.47"
$ws.Range("V33").Formula = "1.72"
$ws.Range("U34").Formula = "1.87"
$ws.Range("V34").Formula = "1.87"
$ws.Range("Q36").Formula = "1.95"
$ws.Range("R36").Formula = "1.8"
$ws.Range("V38").Formula = "1.69"
$ws.Range("U39").Formula = "1.69"
